$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    "H2" = 0.1103252621312025
    "B3" = 0.0704400154605817
    "H3" = 0.1807652775917842
    "B4" = 0.07385558919791943
    "C4" = $null
    "D4" = $null
    "E4" = $null
    "F4" = $null
    "G4" = $null
    "H4" = 0.184180851329122
    "B5" = 0.09775034134038381
    "H5" = 0.2080756034715864
    "B6" = 0.03325985873978371
    "C6" = 0.004473998508572849
    "D6" = 3.424008379840878
    "E6" = 0.04199102364994203
    "F6" = 0.02447910647282289
    "G6" = 0.04204061100674416
    "H6" = 0.1435851208709862
    "B7" = 0.0215173363589086
    "C7" = 0.003392600486410443
    "D7" = 2.282667270780109
    "E7" = 0.01055812936102949
    "F7" = 0.0148560266571679
    "G7" = 0.02817864606064947
    "H7" = 0.1318425984901112
    "B8" = 0.01772736247832942
    "C8" = 0.002272195427779487
    "D8" = 1.406905289337423
    "E8" = 0.007850944809638025
    "F8" = 0.01326402396897121
    "G8" = 0.02219070098768769
    "H8" = 0.128052624609532
    "B9" = -0.004661347441169542
    "C9" = 0.001404839124170956
    "D9" = -1.185296486376566
    "E9" = 0.01305524404181429
    "F9" = -0.007424153017588027
    "G9" = -0.001898541864751114
    "H9" = 0.105663914690033
    "B10" = -0.00368091932506366
    "C10" = 0.00136732187542834
    "D10" = -1.074489543303176
    "E10" = 0.01162710166436559
    "F10" = -0.006365935177336264
    "G10" = -0.0009959034727909342
    "H10" = 0.1066443428061389
    "B11" = 0.02280907856625295
    "H11" = 0.1331343406974555
    "B12" = 0.03855261159719801
    "H12" = 0.1488778737284006
    "B13" = 0.04638141113000906
    "H13" = 0.1567066732612116
    "B14" = 0.05111241853132317
    "C14" = 0.007902941855811025
    "D14" = 11.4731162499843
    "E14" = 0.05140022376806589
    "F14" = 0.03561372598762831
    "G14" = 0.06661111107501799
    "H14" = 0.1614376806625257
    "B15" = 0.05428383737247126
    "H15" = 0.1646090995036738
    "B16" = 0.05753670300054921
    "H16" = 0.1678619651317518
    "B17" = 0.06111087613414442
    "C17" = 0.007994937199944796
    "D17" = 12.6138069044918
    "E17" = 0.03049948613231279
    "F17" = 0.04541438353591733
    "G17" = 0.07680736873237144
    "H17" = 0.171436138265347
    "B18" = -0.1103252621312025
    "C18" = 0.01079776571032865
    "D18" = -17.82924910047596
    "E18" = 0.02540323914022031
    "F18" = -0.131503506623016
    "G18" = -0.08914701763938887
    "B19" = 0.06053595480278799
    "C19" = 0.007702350448359535
    "D19" = 13.22787741156001
    "E19" = 0.03246888672319957
    "F19" = 0.04542557724254849
    "G19" = 0.0756463323630274
    "H19" = 0.1708612169339905
    "B20" = 0.06410739842610239
    "C20" = 0.007714953658747629
    "D20" = 13.16701210081154
    "E20" = 0.03671214015031309
    "F20" = 0.04897271753641649
    "G20" = 0.07924207931578836
    "H20" = 0.1744326605573049
    "B21" = 0.06493580266304178
    "C21" = 0.008091728137715032
    "D21" = 12.96415324187086
    "E21" = 0.04375057291111382
    "F21" = 0.04905205677625769
    "G21" = 0.08081954854982591
    "H21" = 0.1752610647942443
    "B22" = 0.06671571238734235
    "C22" = 0.007376203888435457
    "D22" = 13.04364429162031
    "E22" = 0.05029477987668156
    "F22" = 0.05224472843677717
    "G22" = 0.0811866963379078
    "H22" = 0.1770409745185449
    "B23" = 0.06604278817599865
    "C23" = 0.007459508090423148
    "D23" = 12.99509838570167
    "E23" = 0.05522547490828641
    "F23" = 0.05139904981892887
    "G23" = 0.08068652653306826
    "H23" = 0.1763680503072012
    "B24" = 0.06708940406478939
    "C24" = 0.007363541803498185
    "D24" = 12.48003567790357
    "E24" = 0.0651694139027521
    "F24" = 0.05264844764493008
    "G24" = 0.08153036048464858
    "H24" = 0.1774146661959919
    "B25" = 0.06646527819494623
    "C25" = 0.007271584233597669
    "D25" = 12.21144560552265
    "E25" = 0.05788372980018914
    "F25" = 0.05220423129103739
    "G25" = 0.08072632509885527
    "H25" = 0.1767905403261488
    "B26" = 0.06632014545162199
    "C26" = 0.007410540842542534
    "D26" = 11.90538710729939
    "E26" = 0.0685210476628932
    "F26" = 0.05178683286311346
    "G26" = 0.08085345804013061
    "H26" = 0.1766454075828245
    "B27" = 0.06772988741370799
    "C27" = 0.007407180514326568
    "D27" = 11.81007766279104
    "E27" = 0.06806502373730966
    "F27" = 0.05320223046072899
    "G27" = 0.08225754436668717
    "H27" = 0.1780551495449105
    "B28" = 0.07052531004476637
    "C28" = 0.007581069563438876
    "D28" = 11.72173514496359
    "E28" = 0.08494648414399333
    "F28" = 0.05565808320331801
    "G28" = 0.08539253688621451
    "H28" = 0.1808505721759689
    "B29" = -0.001817690184450246
    "C29" = 0.001277505654093689
    "D29" = -0.8963953259688853
    "E29" = 0.0173566803027261
    "F29" = -0.004325435997626286
    "G29" = 0.000690055628725829
    "H29" = 0.1085075719467523
}

foreach ($ref in $changes.Keys) {
    $ws.Range($ref).Value = $changes[$ref]
}

Write-Host "Applied $($changes.Count) cell updates"
